# "Open Budget category for model.xlsx" — re-split the INC/Interest
# received boundary row (28) and retarget the following "Non-tax" row (29)
# onto the new boundary, and drop the one-off gray-fill formatting that had
# been applied to those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28:29 previously carried a custom (gray-fill) row style; the new
# version of the sheet goes back to the sheet's plain default formatting.
$ws.Rows("28:29").ClearFormats()

# Row 28 (INC / Interest received) now ends at 21050001 instead of
# 21059999 ...
$ws.Range("D28").Value = 21050001

# ... and row 29 (INC / Non-tax) now starts right after it, at 21050002,
# instead of 21060000. D29's upper bound (24109999) is unchanged.
$ws.Range("C29").Value = 21050002

# Reflect the user's new selection/active cell on the sheet.
[void]$ws.Range("C28:D29").Select()
